$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while preserving it as literal Text,
# even when the text looks like a number or date (e.g. "460", "2",
# "82921", "9/9/2022"). Excel's normal Value auto-detects/coerces such
# strings into numbers/dates, so we force the cell to Text format first
# and reset the cell style back to Normal afterwards so no stray
# NumberFormat / style is left behind on the cell.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# New data rows appended to the JOBS sheet (rows 59-60), matching the
# existing row layout: jobNumber, modelNumber, serialNumber, voltage,
# unloaders, statorStatus, incomingNumber, scrap, notes, enteredBy,
# enteredOn, _isDeleted, deletedBy, deletedOn, warranty
$newRows = @(
    @{ Row = 59; JobNumber = 71311; ModelNumber = "O6E3575 661"; SerialNumber = "0920UE9854"; Voltage = "460"; Unloaders = "2"; StatorStatus = "GOOD"; IncomingNumber = "82921"; Scrap = "NO"; Notes = "2 e unl"; EnteredBy = "ravi"; EnteredOn = "9/9/2022"; IsDeleted = $false; DeletedBy = "N/A"; DeletedOn = "N/A"; Warranty = "NO" },
    @{ Row = 60; JobNumber = 71312; ModelNumber = "O6E3575661"; SerialNumber = "5015UE6053"; Voltage = "460"; Unloaders = "2"; StatorStatus = "GOOD"; IncomingNumber = "82923"; Scrap = "NO"; Notes = "2 e unl"; EnteredBy = "ravi"; EnteredOn = "9/9/2022"; IsDeleted = $false; DeletedBy = "N/A"; DeletedOn = "N/A"; Warranty = "NO" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # A: jobNumber - plain number
    $ws.Cells.Item($row, 1).Value = $r.JobNumber

    # B: modelNumber - plain text
    $ws.Cells.Item($row, 2).Value = $r.ModelNumber

    # C: serialNumber - plain text
    $ws.Cells.Item($row, 3).Value = $r.SerialNumber

    # D: voltage - numeric-looking text, keep as Text
    Set-TextValue $ws.Cells.Item($row, 4) $r.Voltage

    # E: unloaders - numeric-looking text, keep as Text
    Set-TextValue $ws.Cells.Item($row, 5) $r.Unloaders

    # F: statorStatus - plain text
    $ws.Cells.Item($row, 6).Value = $r.StatorStatus

    # G: incomingNumber - numeric-looking text, keep as Text
    Set-TextValue $ws.Cells.Item($row, 7) $r.IncomingNumber

    # H: scrap - plain text
    $ws.Cells.Item($row, 8).Value = $r.Scrap

    # I: notes - plain text
    $ws.Cells.Item($row, 9).Value = $r.Notes

    # J: enteredBy - plain text
    $ws.Cells.Item($row, 10).Value = $r.EnteredBy

    # K: enteredOn - date-looking text, keep as Text
    Set-TextValue $ws.Cells.Item($row, 11) $r.EnteredOn

    # L: _isDeleted - boolean
    $ws.Cells.Item($row, 12).Value = $r.IsDeleted

    # M: deletedBy - plain text
    $ws.Cells.Item($row, 13).Value = $r.DeletedBy

    # N: deletedOn - plain text
    $ws.Cells.Item($row, 14).Value = $r.DeletedOn

    # O: warranty - plain text
    $ws.Cells.Item($row, 15).Value = $r.Warranty
}
